$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2 through 27
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22)
$ws.Range("C2:C27").Value = 45221
